# Auto-generated: apply scheduled market-price/profit refresh to Sheets
# (mirrors the per-cell numeric updates captured in the commit diff)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1666736.5
$ws.Range("I6").Value = 1666736.5
$ws.Range("K6").Value = 5000209.5
$ws.Range("M6").Value = -5000097.5
# Row 8
$ws.Range("H8").Value = 902.7646999999999
$ws.Range("I8").Value = 34.7
$ws.Range("K8").Value = 104.1
$ws.Range("M8").Value = 34.89999999999999
# Row 12
$ws.Range("H12").Value = 3180.3333
$ws.Range("I12").Value = 4749.5
$ws.Range("J12").Value = 42
$ws.Range("K12").Value = 4749.5
$ws.Range("L12").Value = 42
$ws.Range("M12").Value = -4579.5
$ws.Range("N12").Value = -382
# Row 17
$ws.Range("H17").Value = 599688.75
$ws.Range("I17").Value = 2021.5238
$ws.Range("J17").Value = 1854790
$ws.Range("K17").Value = 6064.5714
$ws.Range("L17").Value = 5564370
$ws.Range("M17").Value = -5896.5714
$ws.Range("N17").Value = -5564706
# Row 31
$ws.Range("H31").Value = 21.5
$ws.Range("I31").Value = 21.5
$ws.Range("K31").Value = 64.5
$ws.Range("M31").Value = 165.5
# Row 76
$ws.Range("H76").Value = 4481.636
$ws.Range("I76").Value = 3920
$ws.Range("J76").Value = 4949.6665
$ws.Range("K76").Value = 3920
$ws.Range("L76").Value = 4949.6665
$ws.Range("M76").Value = -3605
$ws.Range("N76").Value = -5579.6665
# Row 79
$ws.Range("H79").Value = 4481.636
$ws.Range("I79").Value = 3920
$ws.Range("J79").Value = 4949.6665
$ws.Range("K79").Value = 3920
$ws.Range("L79").Value = 4949.6665
$ws.Range("M79").Value = -2828
$ws.Range("N79").Value = -7133.6665

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7286.675
$ws.Range("I32").Value = 4522.4
$ws.Range("K32").Value = 4522.4
$ws.Range("M32").Value = -4235.4
# Row 61
$ws.Range("H61").Value = 2973.1064
$ws.Range("I61").Value = 2777.9722
$ws.Range("K61").Value = 2777.9722
$ws.Range("M61").Value = -2565.9722
# Row 92
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
# Row 106
$ws.Range("H106").Value = 147500
$ws.Range("J106").Value = 147500
$ws.Range("L106").Value = 147500
$ws.Range("N106").Value = -150024
# Row 132
$ws.Range("H132").Value = 2117.1072
$ws.Range("I132").Value = 1784.0454
$ws.Range("K132").Value = 5352.1362
$ws.Range("M132").Value = -2822.1362
# Row 136
$ws.Range("H136").Value = 2973.1064
$ws.Range("I136").Value = 2777.9722
$ws.Range("K136").Value = 8333.9166
$ws.Range("M136").Value = -5783.9166

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2007.2667
$ws.Range("I94").Value = 1844.1904
$ws.Range("J94").Value = 2387.7778
$ws.Range("K94").Value = 1844.1904
$ws.Range("L94").Value = 2387.7778
$ws.Range("M94").Value = -1393.1904
$ws.Range("N94").Value = -3289.7778
# Row 98
$ws.Range("H98").Value = 25714.285
$ws.Range("J98").Value = 25714.285
$ws.Range("L98").Value = 25714.285
$ws.Range("N98").Value = -31704.285
# Row 99
$ws.Range("H99").Value = 3272
$ws.Range("I99").Value = 1599.8889
$ws.Range("K99").Value = 1599.8889
$ws.Range("M99").Value = -101.8888999999999
# Row 100
$ws.Range("H100").Value = 24285.715
$ws.Range("J100").Value = 24285.715
$ws.Range("L100").Value = 24285.715
$ws.Range("N100").Value = -26449.715
# Row 134
$ws.Range("H134").Value = 1646.5692
$ws.Range("I134").Value = 1653.1968
$ws.Range("K134").Value = 4959.5904
$ws.Range("M134").Value = -2424.5904

$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 287.45456
$ws.Range("I10").Value = 329.8889
$ws.Range("J10").Value = 96.5
$ws.Range("K10").Value = 329.8889
$ws.Range("L10").Value = 96.5
$ws.Range("M10").Value = -190.8889
$ws.Range("N10").Value = -374.5
# Row 22
$ws.Range("H22").Value = 2400
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("N22").Value = -2200
# Row 31
$ws.Range("H31").Value = 20986.16
$ws.Range("I31").Value = 31058
$ws.Range("J31").Value = 4199.7617
$ws.Range("K31").Value = 31058
$ws.Range("L31").Value = 4199.7617
$ws.Range("M31").Value = -30763
$ws.Range("N31").Value = -4789.7617
# Row 34
$ws.Range("H34").Value = 20986.16
$ws.Range("I34").Value = 31058
$ws.Range("J34").Value = 4199.7617
$ws.Range("K34").Value = 31058
$ws.Range("L34").Value = 4199.7617
$ws.Range("M34").Value = -30856
$ws.Range("N34").Value = -4603.7617
# Row 105
$ws.Range("H105").Value = 2090.6316
$ws.Range("I105").Value = 1680.2222
$ws.Range("K105").Value = 1680.2222
$ws.Range("M105").Value = 66.77780000000007
# Row 132
$ws.Range("H132").Value = 2051.5964
$ws.Range("I132").Value = 1621.9318
$ws.Range("K132").Value = 4865.7954
$ws.Range("M132").Value = -2335.7954
# Row 134
$ws.Range("H134").Value = 10387.333
$ws.Range("I134").Value = 3998.1555
$ws.Range("J134").Value = 34346.75
$ws.Range("K134").Value = 11994.4665
$ws.Range("L134").Value = 103040.25
$ws.Range("M134").Value = -9459.466499999999
$ws.Range("N134").Value = -108110.25

$ws = $wb.Worksheets.Item("CUL")
# Row 139
$ws.Range("H139").Value = 1680.8667
$ws.Range("I139").Value = 1309.4166
$ws.Range("K139").Value = 3928.2498
$ws.Range("M139").Value = 1211.7502
# Row 140
$ws.Range("H140").Value = 1418.8948
$ws.Range("I140").Value = 1418.8948
$ws.Range("K140").Value = 4256.6844
$ws.Range("M140").Value = 923.3155999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 504334
$ws.Range("J10").Value = 6501
$ws.Range("L10").Value = 6501
$ws.Range("N10").Value = -6839
# Row 11
$ws.Range("H11").Value = 10200800
# Row 18
$ws.Range("H18").Value = 37103372
$ws.Range("J18").Value = 99500
$ws.Range("L18").Value = 99500
$ws.Range("N18").Value = -100086
# Row 19
$ws.Range("H19").Value = 4666.3335
$ws.Range("I19").Value = 4499.5
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 4499.5
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = -4211.5
$ws.Range("N19").Value = -5576
# Row 126
$ws.Range("H126").Value = 20752.54
$ws.Range("I126").Value = 31098.25
$ws.Range("K126").Value = 93294.75
$ws.Range("M126").Value = -90824.75
# Row 132
$ws.Range("H132").Value = 3247.2188
$ws.Range("I132").Value = 3280
$ws.Range("J132").Value = 2755.5
$ws.Range("K132").Value = 9840
$ws.Range("L132").Value = 8266.5
$ws.Range("M132").Value = -7310
$ws.Range("N132").Value = -13326.5

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4730.933
$ws.Range("I40").Value = 4192.8335
$ws.Range("K40").Value = 4192.8335
$ws.Range("M40").Value = -4056.8335
# Row 122
$ws.Range("H122").Value = 4428.026
$ws.Range("I122").Value = 3840.818
$ws.Range("J122").Value = 7657.6665
$ws.Range("K122").Value = 11522.454
$ws.Range("L122").Value = 22972.9995
$ws.Range("M122").Value = -9072.454000000002
$ws.Range("N122").Value = -27872.9995
# Row 136
$ws.Range("H136").Value = 2645
$ws.Range("I136").Value = 2365.3157
$ws.Range("K136").Value = 7095.9471
$ws.Range("M136").Value = -4545.9471

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1637.4642
$ws.Range("I132").Value = 1637.4642
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4912.392599999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2382.392599999999
$ws.Range("N132").ClearContents()
# Row 136
$ws.Range("H136").Value = 1789.3726
$ws.Range("I136").Value = 1714.9546
$ws.Range("K136").Value = 5144.8638
$ws.Range("M136").Value = -2594.8638

